$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date (column D) and time_last_updated (column E) for rows 2-6
$ws.Range("D2:D6").Value = 44574
$ws.Range("E2:E6").Value = 1642032001

# Update the rates (column F) that changed
$ws.Range("F3").Value = 105.37
$ws.Range("F4").Value = 107.13
$ws.Range("F5").Value = 481.78
